# Insert two new weekly data rows (Primera/Segunda for Provincia de Diguillín,
# fecha 44546) above the existing row 83, shifting the rest of the table down
# by two rows (dimension grows from A1:T196 to A1:T198).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(83).Resize(2).Insert()

# Row 83: "Primera"
$ws.Cells.Item(83, 1).Value = 7
$ws.Cells.Item(83, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(83, 3).Value = "Ñuble"
$ws.Cells.Item(83, 4).Value = 44546
$ws.Cells.Item(83, 5).Value = 16
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100101
$ws.Cells.Item(83, 8).Value = "Berries"
$ws.Cells.Item(83, 9).Value = 100112025
$ws.Cells.Item(83, 10).Value = "Frutilla"
$ws.Cells.Item(83, 11).Value = "Sin especificar"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 120
$ws.Cells.Item(83, 14).Value = 7000
$ws.Cells.Item(83, 15).Value = 7500
$ws.Cells.Item(83, 16).Value = 7250
$ws.Cells.Item(83, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(83, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(83, 19).Value = 1036
$ws.Cells.Item(83, 20).Value = 7

# Row 84: "Segunda"
$ws.Cells.Item(84, 1).Value = 7
$ws.Cells.Item(84, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(84, 3).Value = "Ñuble"
$ws.Cells.Item(84, 4).Value = 44546
$ws.Cells.Item(84, 5).Value = 16
$ws.Cells.Item(84, 6).Value = "Fruta"
$ws.Cells.Item(84, 7).Value = 100101
$ws.Cells.Item(84, 8).Value = "Berries"
$ws.Cells.Item(84, 9).Value = 100112025
$ws.Cells.Item(84, 10).Value = "Frutilla"
$ws.Cells.Item(84, 11).Value = "Sin especificar"
$ws.Cells.Item(84, 12).Value = "Segunda"
$ws.Cells.Item(84, 13).Value = 120
$ws.Cells.Item(84, 14).Value = 6000
$ws.Cells.Item(84, 15).Value = 6500
$ws.Cells.Item(84, 16).Value = 6250
$ws.Cells.Item(84, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(84, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(84, 19).Value = 893
$ws.Cells.Item(84, 20).Value = 7
